# Translations.xlsx - add new translation rows (TranslationsId 29..33) for the
# "first login" / onboarding questionnaire strings.
#
# The shared-strings table is rebuilt at save time from the strings actually
# referenced by cells, in the order those .Value assignments are executed
# (first use wins). The target workbook's shared-strings table lists the ten
# new strings in a specific sequential order that does NOT match the order in
# which they are first read down column D row-by-row (rows 59 and 65 reuse
# strings that logically come "later" in the table). So we first "intern" the
# ten new strings, in the exact desired order, via a scratch cell far outside
# the used range, and only then fill in the real A:D cells (which will simply
# reuse the already-interned shared strings without re-appending/reordering
# them). Finally the scratch cell is cleared so it leaves no trace.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Cells.Item(1000, 1)
$scratch.Value = "On what issues can you help others?"
$scratch.Value = "What issues do you need help with?"
$scratch.Value = "Hangi konularda yardıma ihtiyacınız var?"
$scratch.Value = "Başkalarına hangi konularda yardımcı olabilirsiniz?"
$scratch.Value = "Save"
$scratch.Value = "Kaydol"
$scratch.Value = "It is mandatory to fill this field."
$scratch.Value = "Must contain at least 50 and maximum 300 letters."
$scratch.Value = "Bu alanı doldurmak zorunludur."
$scratch.Value = "En az 50, en fazla 300 karakter içermelidir."
$scratch.ClearContents()

$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = 29
$ws.Cells.Item(58, 3).Value = 1
$ws.Cells.Item(58, 4).Value = "On what issues can you help others?"

$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = 29
$ws.Cells.Item(59, 3).Value = 2
$ws.Cells.Item(59, 4).Value = "Başkalarına hangi konularda yardımcı olabilirsiniz?"

$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = 30
$ws.Cells.Item(60, 3).Value = 1
$ws.Cells.Item(60, 4).Value = "What issues do you need help with?"

$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = 30
$ws.Cells.Item(61, 3).Value = 2
$ws.Cells.Item(61, 4).Value = "Hangi konularda yardıma ihtiyacınız var?"

$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = 31
$ws.Cells.Item(62, 3).Value = 1
$ws.Cells.Item(62, 4).Value = "Save"

$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = 31
$ws.Cells.Item(63, 3).Value = 2
$ws.Cells.Item(63, 4).Value = "Kaydol"

$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = 32
$ws.Cells.Item(64, 3).Value = 1
$ws.Cells.Item(64, 4).Value = "It is mandatory to fill this field."

$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = 32
$ws.Cells.Item(65, 3).Value = 2
$ws.Cells.Item(65, 4).Value = "Bu alanı doldurmak zorunludur."

$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = 33
$ws.Cells.Item(66, 3).Value = 1
$ws.Cells.Item(66, 4).Value = "Must contain at least 50 and maximum 300 letters."

$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = 33
$ws.Cells.Item(67, 3).Value = 2
$ws.Cells.Item(67, 4).Value = "En az 50, en fazla 300 karakter içermelidir."

# Match the saved view state: scrolled so row 46 is at the top, with D67 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1
$null = $ws.Range("D67").Select()
